$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "ECs" sending-cluster rows (old rows 2-4); remaining rows shift up to 2-7
$ws.Rows("2:4").Delete()

# Clear the label columns (A-D) for the remaining data rows so the shared-string table
# can be rebuilt from scratch in the exact order the target workbook uses
$ws.Range("A2:D7").ClearContents()

# Seed the shared-string table in the required order: FAPs, MuSCs, Fgf18, Fgfr2, ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B2").Value = "Fgf18"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"

# Fill in the remaining label cells (reuses the shared strings seeded above)
$ws.Range("A3").Value = "FAPs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A6").Value = "MuSCs"
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B3").Value = "Fgf18"
$ws.Range("B4").Value = "Fgf18"
$ws.Range("B5").Value = "Fgf18"
$ws.Range("B6").Value = "Fgf18"
$ws.Range("B7").Value = "Fgf18"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"

# Fill in the numeric columns E-T
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G2").Value = 8.770835333333332
$ws.Range("G3").Value = 8.770835333333332
$ws.Range("G4").Value = 8.770835333333332
$ws.Range("G5").Value = 0.8194993333333332
$ws.Range("G6").Value = 0.8194993333333332
$ws.Range("G7").Value = 0.8194993333333332
$ws.Range("H2").Value = 26.312506
$ws.Range("H3").Value = 26.312506
$ws.Range("H4").Value = 26.312506
$ws.Range("H5").Value = 2.458498
$ws.Range("H6").Value = 2.458498
$ws.Range("H7").Value = 2.458498
$ws.Range("I2").Value = 0.9145494540267
$ws.Range("I3").Value = 0.9145494540267
$ws.Range("I4").Value = 0.9145494540267
$ws.Range("I5").Value = 0.08545054597330007
$ws.Range("I6").Value = 0.08545054597330007
$ws.Range("I7").Value = 0.08545054597330007
$ws.Range("J2").Value = 0.9145494540267
$ws.Range("J3").Value = 0.9145494540267
$ws.Range("J4").Value = 0.9145494540267
$ws.Range("J5").Value = 0.08545054597330005
$ws.Range("J6").Value = 0.08545054597330005
$ws.Range("J7").Value = 0.08545054597330005
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("M2").Value = 0.036942
$ws.Range("M3").Value = 0.9431116666666667
$ws.Range("M4").Value = 0.7798996666666667
$ws.Range("M5").Value = 0.036942
$ws.Range("M6").Value = 0.9431116666666667
$ws.Range("M7").Value = 0.7798996666666667
$ws.Range("N2").Value = 0.110826
$ws.Range("N3").Value = 2.829335
$ws.Range("N4").Value = 2.339699
$ws.Range("N5").Value = 0.110826
$ws.Range("N6").Value = 2.829335
$ws.Range("N7").Value = 2.339699
$ws.Range("O2").Value = 0.02099032928903418
$ws.Range("O3").Value = 0.5358731102718634
$ws.Range("O4").Value = 0.4431365604391025
$ws.Range("O5").Value = 0.02099032928903418
$ws.Range("O6").Value = 0.5358731102718634
$ws.Range("O7").Value = 0.4431365604391025
$ws.Range("P2").Value = 0.02099032928903418
$ws.Range("P3").Value = 0.5358731102718634
$ws.Range("P4").Value = 0.4431365604391026
$ws.Range("P5").Value = 0.02099032928903418
$ws.Range("P6").Value = 0.5358731102718634
$ws.Range("P7").Value = 0.4431365604391026
$ws.Range("Q2").Value = 0.3240121988839999
$ws.Range("Q3").Value = 8.271877129278888
$ws.Range("Q4").Value = 6.840371552854888
$ws.Range("Q5").Value = 0.03027394437199999
$ws.Range("Q6").Value = 0.7728793820922222
$ws.Range("Q7").Value = 0.6391272569002221
$ws.Range("R2").Value = 2.916109789956
$ws.Range("R3").Value = 74.44689416351
$ws.Range("R4").Value = 61.563343975694
$ws.Range("R5").Value = 0.272465499348
$ws.Range("R6").Value = 6.955914438829999
$ws.Range("R7").Value = 5.752145312101999
$ws.Range("S2").Value = 0.01919669419112685
$ws.Range("S3").Value = 0.4900824604267222
$ws.Range("S4").Value = 0.405270299408851
$ws.Range("S5").Value = 0.001793635097907322
$ws.Range("S6").Value = 0.04579064984514115
$ws.Range("S7").Value = 0.0378662610302516
$ws.Range("T2").Value = 0.01919669419112686
$ws.Range("T3").Value = 0.4900824604267222
$ws.Range("T4").Value = 0.405270299408851
$ws.Range("T5").Value = 0.001793635097907322
$ws.Range("T6").Value = 0.04579064984514115
$ws.Range("T7").Value = 0.0378662610302516
